# daily auto push: 2025-10-03 01:58 UTC
# Append the new daily data row (row 54) to Sheet1, right after the
# existing last row (53), and keep the same "plain text" treatment the
# rest of the sheet uses for the date/weekday columns (so Excel does not
# auto-convert the "2025/10/03" string into a date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 54

# Temporarily force column A to text format so the date-like string is
# stored verbatim, then drop the format override so the cell ends up
# with the sheet's default (unstyled) formatting, matching the other
# data rows.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025/10/03"
$ws.Range("A$newRow").ClearFormats()

$ws.Range("B$newRow").Value = "金"
$ws.Range("C$newRow").Value = 9
$ws.Range("D$newRow").Value = 27
